$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.955.69'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.830.59'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  +0.75%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.30'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4580'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3708'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07188'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8777'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07798'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.66'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.836.17'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.339'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.406'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.26'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -5.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('D17').ClearFormats()
$ws.Range('E18').Value = '  -1.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.985.64'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.52'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.012'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.058.02'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.43'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.029'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +7.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.54'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.966'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.13'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.939'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08801'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.034'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7517'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.486'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.14%  '
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.563'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.093'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01948'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05158'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.894'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.964'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4990'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1603'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.324'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4691'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.92%  '
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.54'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.614'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06120'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.62'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.48%  '
